# Regenerate merged AHB files
# - Rename the "_old"/"_new" suffixed header labels to "_FV2310"/"_FV2404"
# - Turn the data range into an Excel Table ("Table1")
# - Freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row (row 1) cell values -----------------------------
# Columns A..J used to end in "_old" -> now end in "_FV2310"
$headersNew = @("Segmentname_FV2310","Segmentgruppe_FV2310","Segment_FV2310","Datenelement_FV2310","Segment ID_FV2310","Code_FV2310","Qualifier_FV2310","Beschreibung_FV2310","Bedingungsausdruck_FV2310","Bedingung_FV2310")

# Columns L..U used to end in "_new" -> now end in "_FV2404"
$headersNew2 = @("Segmentname_FV2404","Segmentgruppe_FV2404","Segment_FV2404","Datenelement_FV2404","Segment ID_FV2404","Code_FV2404","Qualifier_FV2404","Beschreibung_FV2404","Bedingungsausdruck_FV2404","Bedingung_FV2404")

# Columns A..J (1..10)
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headersNew[$i]
}

# Column K (11) = "diff" stays the same

# Columns L..U (12..21)
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $headersNew2[$i]
}

# --- 2. Convert the range into a table -------------------------------------
$tableRange = $ws.Range("A1:U72")
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"
$lo.TableStyle = ""

# --- 3. Freeze the header row -----------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Host "done"
